# transactions.xlsx - add "Discounted Price" column and chart it (3x)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header tweaks -------------------------------------------------------
$ws.Range("C1").Value = "Price"
$ws.Range("D1").Value = "Discounted Price"

# --- New discounted-price data in column D --------------------------------
$ws.Range("D2").Value = 5.36
$ws.Range("D3").Value = 6.25
$ws.Range("D4").Value = 7.16

# --- Drop the old, unlabeled data that used to live in column E -----------
$ws.Range("E2:E4").ClearContents()

# --- Add three identical bar charts plotting the discounted price ---------
for ($i = 1; $i -le 3; $i++) {
    $co = $ws.ChartObjects().Add(304, 20, 300, 150)
    $co.Chart.ChartType = 51
    $ser = $co.Chart.SeriesCollection().NewSeries()
    $ser.Values = "='Sheet1'!`$D`$2:`$D`$4"
}
